$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0161
$ws.Range("E2").Value = 0.01215
$ws.Range("F2").Value = 0.1654
$ws.Range("G2").Value = 0.09703553524722451
$ws.Range("H2").Value = 0.09703553524722451
$ws.Range("I2").Value = 0.08972305340308846
$ws.Range("J2").Value = 0.07484730222079432
$ws.Range("K2").Value = 8976.74
$ws.Range("L2").Value = 0.04929466013046412
$ws.Range("M2").Value = 465.97
$ws.Range("N2").Value = 0.004142791478813947
$ws.Range("O2").Value = 0.05190859933561627
$ws.Range("P2").Value = 465.97
$ws.Range("Q2").Value = 0.004142791478813947
$ws.Range("R2").Value = 0.05190859933561627
$ws.Range("U2").Value = 30655.85
$ws.Range("V2").Value = 0.2725514392681901
$ws.Range("W2").Value = 0.09889207134327588
$ws.Range("X2").Value = 0.0410960917177257
$ws.Range("Y2").Value = 0.05779597962555018
$ws.Range("Z2").Value = 1.585696489752046
$ws.Range("AA2").Value = 0.0805393771442103
$ws.Range("AB2").Value = 0.04035682311219418
$ws.Range("AC2").Value = 0.04018255403201612
$ws.Range("AD2").Value = 51565.64
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 51565.64
$ws.Range("AG2").Value = 20909.79
$ws.Range("AH2").Value = 0.3143423301240517
$ws.Range("AI2").Value = 0.3146999198201061
$ws.Range("AJ2").Value = 0.1567602231970125
$ws.Range("AK2").Value = 0.1569795922228331
$ws.Range("AL2").Value = 1608.67
$ws.Range("AM2").Value = 1608.67
$ws.Range("AN2").Value = 2.877228966888703
$ws.Range("AO2").Value = 10.15677547290619
$ws.Range("AP2").Value = 1.166712048557135
$ws.Range("AQ2").Value = 10.15677547290619

# Row 3
$ws.Range("D3").Value = 0.0554
$ws.Range("E3").Value = -0.00726
$ws.Range("F3").Value = 0.265
$ws.Range("G3").Value = 0.06947996956422679
$ws.Range("H3").Value = 0.06947996956422679
$ws.Range("I3").Value = 0.05137541322580506
$ws.Range("J3").Value = 0.0393312441069017
$ws.Range("K3").Value = 823.2
$ws.Range("L3").Value = 0.01769401730712189
$ws.Range("M3").Value = 418
$ws.Range("N3").Value = 0.04258093433571705
$ws.Range("O3").Value = 0.5077745383867832
$ws.Range("P3").Value = 418
$ws.Range("Q3").Value = 0.04258093433571705
$ws.Range("R3").Value = 0.5077745383867832
$ws.Range("U3").Value = 4533.8
$ws.Range("V3").Value = 0.4618503351465884
$ws.Range("W3").Value = 0.07378393639810342
$ws.Range("X3").Value = 0.05557382269949658
$ws.Range("Y3").Value = 0.01821011369860684
$ws.Range("Z3").Value = 4.237910020859712
$ws.Range("AA3").Value = 0.1666822735335182
$ws.Range("AB3").Value = 0.03858882354160514
$ws.Range("AC3").Value = 0.1280934499919131
$ws.Range("AD3").Value = 6777.8
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 6777.8
$ws.Range("AG3").Value = 2244
$ws.Range("AH3").Value = 0.4084389914670009
$ws.Range("AI3").Value = 0.2554874118414258
$ws.Range("AJ3").Value = 0.1860603950052236
$ws.Range("AK3").Value = 0.1020227232429041
$ws.Range("AL3").Value = 439.2
$ws.Range("AM3").Value = 439.2
$ws.Range("AN3").Value = 2.520002974419988
$ws.Range("AO3").Value = 5.442167577413478
$ws.Range("AP3").Value = 0.8343248066627008
$ws.Range("AQ3").Value = 5.442167577413478

# Row 4
$ws.Range("D4").Value = 0.0182
$ws.Range("E4").Value = 0.0131
$ws.Range("F4").Value = 0.0658
$ws.Range("G4").Value = 0.1108851788068531
$ws.Range("H4").Value = 0.1108851788068531
$ws.Range("I4").Value = 0.106419978159594
$ws.Range("J4").Value = 0.07895207973573588
$ws.Range("K4").Value = 8041.3
$ws.Range("L4").Value = 0.06149343717212626
$ws.Range("U4").Value = 25442.7
$ws.Range("V4").Value = 0.2517244859191718
$ws.Range("W4").Value = 0.09889207134327588
$ws.Range("X4").Value = 0.05016652768895191
$ws.Range("Y4").Value = 0.04872554365432397
$ws.Range("Z4").Value = 1.267438242737276
$ws.Range("AA4").Value = 0.1000668852007144
$ws.Range("AB4").Value = 0.03860389418259323
$ws.Range("AC4").Value = 0.06146299101812117
$ws.Range("AD4").Value = 44758.1
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 44758.1
$ws.Range("AG4").Value = 19315.4
$ws.Range("AH4").Value = 0.3069161231748652
$ws.Range("AI4").Value = 0.3284899316791201
$ws.Range("AJ4").Value = 0.1604415685818472
$ws.Range("AK4").Value = 0.1743087805045329
$ws.Range("AL4").Value = 1168.4
$ws.Range("AM4").Value = 1168.4
$ws.Range("AN4").Value = 2.946511566667983
$ws.Range("AO4").Value = 11.91047586442999
$ws.Range("AP4").Value = 1.271569827915366
$ws.Range("AQ4").Value = 11.91047586442999

# Row 5
$ws.Range("D5").Value = -0.0766
$ws.Range("E5").Value = 0.0112
$ws.Range("G5").Value = -0.01118346545866365
$ws.Range("H5").Value = -0.01118346545866365
$ws.Range("I5").Value = 0.01380237825594564
$ws.Range("J5").Value = 0.01204290195803009
$ws.Range("K5").Value = 99.4
$ws.Range("L5").Value = 0.02345224613061533
$ws.Range("M5").Value = 42.8
$ws.Range("N5").Value = 0.04228413357044062
$ws.Range("O5").Value = 0.4305835010060362
$ws.Range("P5").Value = 42.8
$ws.Range("Q5").Value = 0.04228413357044062
$ws.Range("R5").Value = 0.4305835010060362
$ws.Range("U5").Value = 675.9
$ws.Range("V5").Value = 0.6677534084173088
$ws.Range("W5").Value = 0.1092788038698329
$ws.Range("X5").Value = 0.0410960917177257
$ws.Range("Y5").Value = 0.0681827121521072
$ws.Range("Z5").Value = 6.687705124968442
$ws.Range("AA5").Value = 0.0805393771442103
$ws.Range("AB5").Value = 0.04035682311219418
$ws.Range("AC5").Value = 0.04018255403201612
$ws.Range("AD5").Value = 27.8
$ws.Range("AF5").Value = 27.8
$ws.Range("AG5").Value = -648.1
$ws.Range("AH5").Value = 0.02673076923076923
$ws.Range("AI5").Value = 0.02765068629401234
$ws.Range("AJ5").Value = -1.78000549299643
$ws.Range("AK5").Value = -1.966919575113809
$ws.Range("AN5").Value = 0.4648829431438127
$ws.Range("AP5").Value = -10.83779264214047

# Row 6
$ws.Range("G6").Value = -0.05935582822085889
$ws.Range("H6").Value = -0.05935582822085889
$ws.Range("I6").Value = -0.1150306748466258
$ws.Range("J6").Value = -0.1150306748466258
$ws.Range("K6").Value = -3.86
$ws.Range("L6").Value = -0.05920245398773006
$ws.Range("U6").Value = 3.45
$ws.Range("V6").Value = 0.008443465491923641
$ws.Range("W6").Value = -0.05237449118046133
$ws.Range("X6").Value = 0.04060001064693675
$ws.Range("Y6").Value = -0.09297450182739808
$ws.Range("Z6").Value = 1.174774774774775
$ws.Range("AA6").Value = -0.1351351351351351
$ws.Range("AB6").Value = 0.04047166646450615
$ws.Range("AC6").Value = -0.1756068015996413
$ws.Range("AD6").Value = 1.94
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 1.94
$ws.Range("AG6").Value = -1.51
$ws.Range("AH6").Value = 0.004725483509524041
$ws.Range("AI6").Value = 0.02847079542119166
$ws.Range("AJ6").Value = -0.003709253482031001
$ws.Range("AK6").Value = -0.02334209305920545
$ws.Range("AL6").Value = 1.07
$ws.Range("AM6").Value = 1.07
$ws.Range("AN6").Value = -0.2844574780058651
$ws.Range("AO6").Value = -7.009345794392523
$ws.Range("AP6").Value = 0.2214076246334311
$ws.Range("AQ6").Value = -7.009345794392523

# Row 7
$ws.Range("D7").Value = 0.014
$ws.Range("E7").Value = 0.124
$ws.Range("G7").Value = -0.02121390689451974
$ws.Range("H7").Value = -0.02121390689451974
$ws.Range("I7").Value = -0.03633863681005696
$ws.Range("J7").Value = -0.02874523298606864
$ws.Range("K7").Value = 16.7
$ws.Range("L7").Value = 0.03280298566097033
$ws.Range("M7").Value = 5.17
$ws.Range("N7").Value = 0.03108839446782922
$ws.Range("O7").Value = 0.3095808383233533
$ws.Range("P7").Value = 5.17
$ws.Range("Q7").Value = 0.03108839446782922
$ws.Range("R7").Value = 0.3095808383233533
$ws.Range("X7").Value = 0.04049632828323475
$ws.Range("AB7").Value = 0.04049632828323475
